# Apply updated "想去人数" (F column) counts across the four sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2754
$ws1.Range("F4").Value = 1086
$ws1.Range("F5").Value = 20042
$ws1.Range("F6").Value = 82
$ws1.Range("F7").Value = 2324
$ws1.Range("F8").Value = 761
$ws1.Range("F10").Value = 455
$ws1.Range("F11").Value = 703
$ws1.Range("F12").Value = 252
$ws1.Range("F14").Value = 73
$ws1.Range("F15").Value = 384
$ws1.Range("F16").Value = 86
$ws1.Range("F17").Value = 277
$ws1.Range("F18").Value = 174
$ws1.Range("F19").Value = 219
$ws1.Range("F22").Value = 105

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 198
$ws2.Range("F5").Value = 21
$ws2.Range("F7").Value = 297
$ws2.Range("F10").Value = 16
$ws2.Range("F13").Value = 1
$ws2.Range("F16").Value = 102

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6022
$ws3.Range("F4").Value = 599

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6022
$ws4.Range("F4").Value = 599
$ws4.Range("F5").Value = 198
$ws4.Range("F8").Value = 2754
$ws4.Range("F9").Value = 1086
$ws4.Range("F10").Value = 20042
$ws4.Range("F12").Value = 21
$ws4.Range("F13").Value = 82
$ws4.Range("F15").Value = 297
$ws4.Range("F16").Value = 2324
$ws4.Range("F17").Value = 761
$ws4.Range("F20").Value = 455
$ws4.Range("F21").Value = 703
$ws4.Range("F22").Value = 252
$ws4.Range("F25").Value = 73
$ws4.Range("F26").Value = 16
$ws4.Range("F28").Value = 384
$ws4.Range("F29").Value = 86
$ws4.Range("F31").Value = 1
$ws4.Range("F32").Value = 277
$ws4.Range("F34").Value = 174
$ws4.Range("F36").Value = 219
$ws4.Range("F37").Value = 102
$ws4.Range("F38").Value = 102
$ws4.Range("F49").Value = 105
